$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1

# Swap the contents of columns D (codeforiati:group-name) and E
# (codeforiati:group-code) for every row, header included. The source
# data had these two columns reversed; this brings column D in line
# with "group-code" values and column E with "group-name" values.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
}
